# Apply the edit described by the diff:
# - Insert two new rows at position 254 (shifting existing rows 254-336 down to 256-338)
# - Populate the two new rows (254 and 255) with new data

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows before the current row 254 (shifts rows 254:336 down to 256:338)
$ws.Rows("254:255").Insert()

# --- New row 254 ---
$ws.Range("A254").Value = 9
$ws.Range("B254").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C254").Value = "Metropolitana"
$ws.Range("D254").Value = 44627
$ws.Range("E254").Value = 13
$ws.Range("F254").Value = 100112052
$ws.Range("G254").Value = "Albahaca"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Primera"
$ws.Range("J254").Value = 79
$ws.Range("K254").Value = 3000
$ws.Range("L254").Value = 3500
$ws.Range("M254").Value = 3253
$ws.Range("N254").Value = '$/docena de matas'
$ws.Range("O254").Value = "Región Metropolitana"
$ws.Range("P254").Value = 542
$ws.Range("Q254").Value = 6
$ws.Range("R254").Value = "Hortaliza"

# --- New row 255 ---
$ws.Range("A255").Value = 9
$ws.Range("B255").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C255").Value = "Metropolitana"
$ws.Range("D255").Value = 44627
$ws.Range("E255").Value = 13
$ws.Range("F255").Value = 100112052
$ws.Range("G255").Value = "Albahaca"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Segunda"
$ws.Range("J255").Value = 43
$ws.Range("K255").Value = 2500
$ws.Range("L255").Value = 2500
$ws.Range("M255").Value = 2500
$ws.Range("N255").Value = '$/docena de matas'
$ws.Range("O255").Value = "Región Metropolitana"
$ws.Range("P255").Value = 417
$ws.Range("Q255").Value = 6
$ws.Range("R255").Value = "Hortaliza"
